$wb = $excel.ActiveWorkbook

# Productdata sheet: G2 changes from 40 to 70
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("G2").Value = 70

# Preserve the empty (blank-string) H2:H11 cells, which the engine would
# otherwise resolve to shared-string index 0 ("Name") on save if left untouched.
$wsProductdata.Range("H2:H11").Value = ""

# ForecastedAverageDemand sheet: B9:B11 change from 0 to 100
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

# ForcastedStandardDeviation sheet: B9:B11 change from 0 to new values
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775
